$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.729.37'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.79'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.38'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5318'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3727'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07165'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.43'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8862'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08185'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.868.87'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +27.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.18'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.286'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.80'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008476'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.769.54'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.967'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.361'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.289'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.73'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.85%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.730'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.683'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.615'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09113'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7999'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05001'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.168'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.943'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6079'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.661'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.174'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01939'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.062'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.479'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5190'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.708'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '114.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1490'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.633'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.926'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.35'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06060'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.06'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.78%  '
